# Atualização de bases das ligas, do dia: 04-04-2024 às 23:22
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Swap the full (B:AC) content between row pairs that got re-ordered.
#    Column A ("id") is left untouched - it always equals (row number - 2)
#    and is already correct on both sides of every swap.
# ---------------------------------------------------------------------------
function Swap-RowContent($rowA, $rowB) {
    $rangeA = $ws.Range("B$rowA`:AC$rowA")
    $rangeB = $ws.Range("B$rowB`:AC$rowB")
    $valsA = $rangeA.Value2()
    $valsB = $rangeB.Value2()
    $rangeA.Value2 = $valsB
    $rangeB.Value2 = $valsA
}

Swap-RowContent 108 110
Swap-RowContent 112 114
Swap-RowContent 137 138
Swap-RowContent 139 140

# ---------------------------------------------------------------------------
# 2) Two new fixtures were added (rows 170 and 175 in the final sheet),
#    shifting the old rows 170-174 down (old 174 ends up at 176).
#    Capture the old content (B:AC) of rows 170-174 first ...
# ---------------------------------------------------------------------------
$old170 = $ws.Range("B170:AC170").Value2()
$old171 = $ws.Range("B171:AC171").Value2()
$old172 = $ws.Range("B172:AC172").Value2()
$old173 = $ws.Range("B173:AC173").Value2()
$old174 = $ws.Range("B174:AC174").Value2()

# ... then prepare the two new sheet rows (175 and 176) with the right
# border/centre formatting used by every data row (copy format from the
# last existing data row, which already has the same sparse column layout).
$ws.Range("A174").Copy()
$ws.Range("A175:A176").PasteSpecial(-4122)
$ws.Range("E174").Copy()
$ws.Range("E175:E176").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# row "id" (col A) is always (row number - 2) - keep it consistent on the
# two brand-new physical rows created by the shift.
$ws.Range("A176").Value2 = 174

# ... and write the shifted content into its new home (bottom-up so nothing
# gets clobbered before it has been read - values were cached above anyway).
$ws.Range("B176:AC176").Value2 = $old174
$ws.Range("B174:AC174").Value2 = $old173
$ws.Range("B173:AC173").Value2 = $old172
$ws.Range("B172:AC172").Value2 = $old171
$ws.Range("B171:AC171").Value2 = $old170

# ---------------------------------------------------------------------------
# 3) Fill in the two brand-new fixtures (row 170 and row 175).
# ---------------------------------------------------------------------------
$ws.Range("A170").Value2 = 168
$ws.Range("B170").Value2 = 7723553
$ws.Range("C170").Value2 = "Chile Primera Division"
$ws.Range("D170").Value2 = "Chile Primera Division"
$ws.Range("E170").Value2 = 45388.52083333334
$ws.Range("F170").Value2 = "Deportes Copiapo"
$ws.Range("G170").Value2 = "Deportes Iquique"
$ws.Range("K170").Value2 = 3.2
$ws.Range("L170").Value2 = 3.4
$ws.Range("M170").Value2 = 2.2
$ws.Range("N170").Value2 = 3.2
$ws.Range("O170").Value2 = 3.6
$ws.Range("P170").Value2 = 2.15
$ws.Range("Q170").Value2 = 0.25
$ws.Range("R170").Value2 = 1.975
$ws.Range("S170").Value2 = 1.875
$ws.Range("T170").Value2 = 2.75
$ws.Range("U170").Value2 = 2.05
$ws.Range("V170").Value2 = 1.8
$ws.Range("W170").Value2 = 0
$ws.Range("X170").Value2 = 0
$ws.Range("Y170").Value2 = 0
$ws.Range("Z170").Value2 = 0
$ws.Range("AA170").Value2 = 0

$ws.Range("A175").Value2 = 173
$ws.Range("B175").Value2 = 7723552
$ws.Range("C175").Value2 = "Chile Primera Division"
$ws.Range("D175").Value2 = "Chile Primera Division"
$ws.Range("E175").Value2 = 45389.77083333334
$ws.Range("F175").Value2 = "Coquimbo Unido"
$ws.Range("G175").Value2 = "Cobreloa"
$ws.Range("K175").Value2 = 2.4
$ws.Range("L175").Value2 = 3.4
$ws.Range("M175").Value2 = 2.8
$ws.Range("N175").Value2 = 2.3
$ws.Range("O175").Value2 = 3.4
$ws.Range("P175").Value2 = 2.9
$ws.Range("Q175").Value2 = -0.25
$ws.Range("R175").Value2 = 2.025
$ws.Range("S175").Value2 = 1.825
$ws.Range("T175").Value2 = 2.5
$ws.Range("U175").Value2 = 2
$ws.Range("V175").Value2 = 1.85
$ws.Range("W175").Value2 = 0
$ws.Range("X175").Value2 = 0
$ws.Range("Y175").Value2 = 0
$ws.Range("Z175").Value2 = 0
$ws.Range("AA175").Value2 = 0

Write-Host "Done applying Chile Primera Division update"
